$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates to whole numbers
$ws.Range("Q2").Value = 630224
$ws.Range("R2").Value = 6534431

# Remove the time values (Starttid / Sluttid) entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
